$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Shared-string bookkeeping note: G7 originally carries the *old* "Create"
# steps text. The new content re-shuffles which cell owns which text, so we
# clear G7 first (releasing that slot) before re-populating things in the
# same order the workbook's string table ends up using them.
# ---------------------------------------------------------------------------
$ws.Range("G7").Value = ""

# ---------------------------------------------------------------------------
# Row 8 ("Modify" test case): Test Data / Steps
# ---------------------------------------------------------------------------
$ws.Range("G8").Value = @"
1.Log in to KanbanBoard application
2.Hover over a board 
3.Click on 'Modify' icon
4.Enter new name of the board and click on 'Cancel' button
5.Hover over a board and click on 'Modify' icon
6.Enter new name of the board and click 'Modify'
7.Hover over a board and click on 'Modify' icon
8.Enter new name of the board and click 'Enter' key
"@

$ws.Range("I8").Value = @"
2.'Modify', 'Remove' icons are displayed
3.Title of the 'board' become editable
4.Title of the 'board' hasn't changed
6.Title of the 'board' has changed
8.Title of the 'board' has changed

"@

# ---------------------------------------------------------------------------
# Row 9 ("Remove/Delete" test case): Test Data / Steps
# ---------------------------------------------------------------------------
$ws.Range("G9").Value = @"
1.Log in to KanbanBoard application
2.Hover over a board 
3.Click on 'Remove' icon
4.In confirmation dialog click on 'No' button
5.Hover over a board and click on 'Remove' icon
6.In confirmation dialog click on 'Yes' button

"@

$ws.Range("I9").Value = @"
2.'Modify', 'Remove' icons are displayed
3.Confirmation dialog appears
4.Board has not been removed
5.Board has been removed

"@

# ---------------------------------------------------------------------------
# Row 7 ("Create" test case): Test Data / Steps (rewritten, simplified flow)
# ---------------------------------------------------------------------------
$ws.Range("G7").Value = @"
1.Log in to KanbanBoard application
2.Click on 'Create' button
3.Enter name of new board into 'title' field
4.Click on 'Cancel' button
5.Click on 'Create' button and enter name of new board
6.Click on 'Create' button
"@

$ws.Range("I7").Value = @"
1.'Create' button is displayed and is active
2.New dialog is displayed
4.Board is not created and not is displayed
6.Board is created and displayed in the board's list
"@

# ---------------------------------------------------------------------------
# Expected Results column (K) for rows 7-9: a "verify board in DB" reminder,
# highlighted with the built-in "Input" cell style (orange fill).
# ---------------------------------------------------------------------------
$ws.Range("K7").Value = "verify board in DB"
$ws.Range("K8").Value = "verify board in DB"
$ws.Range("K9").Value = "verify board in DB"

$ws.Range("K7:K9").Style = "Input"
$ws.Range("K7:K9").WrapText = $true
$ws.Range("K7:K9").VerticalAlignment = -4160

# ---------------------------------------------------------------------------
# New rows 10 & 11 pick up a Priority + sequential ID
# ---------------------------------------------------------------------------
$ws.Range("A10").Value = "High"
$ws.Range("C10").Value = 9

$ws.Range("A11").Value = "High"
$ws.Range("C11").Value = 10

# ---------------------------------------------------------------------------
# Row heights (grew to fit the new text)
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).RowHeight = 89.25
$ws.Rows.Item(8).RowHeight = 165.75
$ws.Rows.Item(9).RowHeight = 102

# ---------------------------------------------------------------------------
# Sheet view: scroll down a bit and move the selection
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select()
